# Publish terminology IG 2.0.0
#
# 1. Update Metadata: Version 1.8.1 -> 1.0.2, Date 2023-01-06 -> 2025-09-22
# 2. Rename existing "Concepts" sheet (which holds the old concept/code
#    table) to "Properties", trim it down to the 3-row FHIR properties
#    table (Code/Uri/Description/Type, status, effectiveDate).
# 3. Add a brand-new "Concepts" sheet (placed after "Properties") that
#    holds the original Level/Code/Display/Definition concept table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet edits
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value = "1.0.2"

# The Date cell ("2023-01-06" -> "2025-09-22") looks like an ISO date,
# so a plain .Value assignment gets auto-converted into a date serial
# number by Excel. Write it as a text-formula first (forces a string
# result), then flatten the formula down to a literal value via
# copy/paste-values so it lands back as plain shared-string text under
# the cell's original (unchanged) style.
$dateCell = $meta.Cells.Item(8, 2)
$dateCell.Formula = '="2025-09-22"'
$dateCell.Copy() | Out-Null
$dateCell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Rename "Concepts" -> "Properties" and replace its data
# ---------------------------------------------------------------------
$propsWs = $wb.Worksheets.Item("Concepts")
$propsWs.Name = "Properties"

# Drop the old concept rows (2..26 data + nothing else) beyond row 3 so
# the sheet's used range becomes A1:D3.
$propsWs.Range("A4:D26").EntireRow.Delete()

$propsData = @(
    @("Code", "Uri", "Description", "Type"),
    @("status", "http://hl7.org/fhir/concept-properties#status", "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired.", "code"),
    @("effectiveDate", "http://hl7.org/fhir/concept-properties#effectiveDate", "The date at which the concept status was last changed.", "dateTime")
)
for ($i = 0; $i -lt $propsData.Count; $i++) {
    $row = $propsData[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $propsWs.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}

# ---------------------------------------------------------------------
# 3. New "Concepts" sheet with the original concept table
# ---------------------------------------------------------------------
$conceptsWs = $wb.Worksheets.Add($null, $propsWs)
$conceptsWs.Name = "Concepts"

$conceptsHeader = @("Level", "Code", "Display", "Definition")
$conceptsData = @(
    @("1", "new-message", "A new message has been sent", $null),
    @("1", "forward-message", "A message has been forwarded and potentially changed", $null),
    @("1", "reply-message", "A message has been replied to", $null),
    @("1", "retract-message", "A message has been retracted/cancelled", $null),
    @("1", "modified-message", "A message has been modified", $null),
    @("1", "carboncopy-message", "An exact copy of a message to a primary receiver has been sent to a carbon copy receiver", $null),
    @("1", "admit-emergency", "Start hospital stay - acute ambulant", $null),
    @("1", "revise-admit-emergency", "Update Start hospital stay - acute ambulant", $null),
    @("1", "cancel-admit-emergency", "Cancellation Start hospital stay - acute ambulant", $null),
    @("1", "admit-inpatient", "Start hospital stay - admitted", $null),
    @("1", "revise-admit-inpatient", "Update Start hospital stay - admitted", $null),
    @("1", "cancel-admit-inpatient", "Cancellation Start hospital stay - admitted", $null),
    @("1", "start-leave-inpatient", "Start leave", $null),
    @("1", "revise-start-leave-inpatient", "Update Start leave", $null),
    @("1", "cancel-start-leave-inpatient", "Cancellation Start leave", $null),
    @("1", "end-leave-inpatient", "End leave", $null),
    @("1", "revise-end-leave-inpatient", "Update End leave", $null),
    @("1", "cancel-end-leave-inpatient", "Cancellation End leave", $null),
    @("1", "discharge-emergency-home", "End hospital stay - patient completed to home/primary sector", $null),
    @("1", "discharge-inpatient-home", "End hospital stay - patient completed to home/primary sector", $null),
    @("1", "revise-discharge-emergency-home", "Update End hospital stay - patient completed to home/primary sector", $null),
    @("1", "revise-discharge-inpatient-home", "Update End hospital stay - patient completed to home/primary sector", $null),
    @("1", "cancel-discharge-emergency-home", "Cancellation End hospital stay - patient completed to home/primary sector", $null),
    @("1", "cancel-discharge-inpatient-home", "Cancellation End hospital stay - patient completed to home/primary sector", $null),
    @("1", "acknowledgement", "Acknowledgement message has been sent", $null)
)

for ($j = 0; $j -lt $conceptsHeader.Count; $j++) {
    $conceptsWs.Cells.Item(1, $j + 1).Value = $conceptsHeader[$j]
}
for ($i = 0; $i -lt $conceptsData.Count; $i++) {
    $row = $conceptsData[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        if ($null -ne $row[$j]) {
            $conceptsWs.Cells.Item($i + 2, $j + 1).Value = $row[$j]
        }
    }
}

# Column A ("Level") holds the numeric-looking text "1" for every data
# row; a plain .Value assignment gets silently coerced into a real
# number by Excel. Force it back to literal shared-string text the same
# way as the Date cell above: write as a text formula, then flatten via
# copy/paste-values across the whole column in one shot.
$levelRange = $conceptsWs.Range("A2:A$($conceptsData.Count + 1)")
foreach ($cell in $levelRange) {
    $cell.Formula = '="1"'
}
$levelRange.Copy() | Out-Null
$levelRange.PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = 0

# Copy the header/data row formatting (fill, border, bold, wrap …) from
# the Properties sheet (which still carries the original "Concepts"
# styles: row 1 = bold header style, rows below = plain bordered style)
# onto the new Concepts sheet.
$propsWs.Range("A1:D1").Copy() | Out-Null
$conceptsWs.Range("A1:D26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$propsWs.Range("A2:D2").Copy() | Out-Null
$conceptsWs.Range("A2:D26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
